$d = $word.ActiveDocument

# Update the date paragraph
$d.Content.Find.Execute("2025-05-24 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-25 Sunday", 2) | Out-Null

# Update the table cells (20 rows x 5 columns), in document order
$t = $d.Tables.Item(1)
$values = @(
    "84-31=53",
    "19+33=52",
    "85-25=60",
    "55+12=67",
    "16+29=45",
    "89-56=33",
    "91-3=88",
    "12+45=57",
    "81+2=83",
    "32+23=55",
    "7+66=73",
    "88-59=29",
    "58-46=12",
    "25-18=7",
    "48+28=76",
    "41+9=50",
    "2+23=25",
    "24-11=13",
    "25+50=75",
    "24+60=84",
    "17+46=63",
    "40+53=93",
    "8+14=22",
    "89-77=12",
    "82+3=85",
    "93-81=12",
    "93-82=11",
    "90-34=56",
    "90-38=52",
    "4+40=44",
    "21+74=95",
    "46+34=80",
    "92-53=39",
    "69-21=48",
    "89-49=40",
    "62+10=72",
    "28+2=30",
    "38+31=69",
    "75+15=90",
    "4+95=99",
    "10+61=71",
    "37-1=36",
    "78-6=72",
    "42-1=41",
    "27+48=75",
    "18+60=78",
    "63-42=21",
    "90-56=34",
    "85-21=64",
    "87+1=88",
    "23+35=58",
    "31+11=42",
    "95-61=34",
    "81-50=31",
    "22+26=48",
    "70-53=17",
    "55+16=71",
    "90-88=2",
    "95-94=1",
    "11+69=80",
    "34+4=38",
    "64-59=5",
    "48+47=95",
    "7+4=11",
    "55+9=64",
    "21-6=15",
    "70+14=84",
    "92-12=80",
    "27+33=60",
    "54-45=9",
    "22+58=80",
    "28+19=47",
    "95-87=8",
    "48-31=17",
    "58-36=22",
    "34+55=89",
    "61+4=65",
    "55-5=50",
    "91-10=81",
    "97-23=74",
    "1+83=84",
    "73-64=9",
    "93-70=23",
    "42-1=41",
    "97-76=21",
    "86-49=37",
    "23+15=38",
    "17+42=59",
    "37-25=12",
    "87+11=98",
    "95-10=85",
    "7+68=75",
    "37+42=79",
    "54+0=54",
    "87-9=78",
    "38+54=92",
    "83-70=13",
    "79+2=81",
    "77-72=5",
    "57-14=43"
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}
